# "changed the whole search flow"
# GBT.xlsx test-data sheet: re-worked the SLA/Activity search rows —
# row 2's Param1 flips from BRO to aff, the stale Actual/Pass-Fail columns
# (R/S) are dropped everywhere, row 3 gains a Q "true" flag, and three new
# rows (6-8) are appended covering the BRO / xyz / TestSLA4 search cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: Param1 BRO -> aff; Q2's quote-prefixed "true" switches from the
# text-numfmt style (s=1) to the general+quotePrefix style (s=3, same one
# row 4 already uses); the leftover Actual/Pass-Fail (R2:S2) go away.
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "aff"
$ws.Range("Q4").Copy()
$ws.Range("Q2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R2:S2").ClearContents()

# ---------------------------------------------------------------------
# Row 3: just "Activity" before; now also carries a Q "true" flag (s=3).
# ---------------------------------------------------------------------
$ws.Range("Q4").Copy()
$ws.Range("Q3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q3").Value = "'true"

# ---------------------------------------------------------------------
# Row 4: drop the leftover Actual/Pass-Fail columns (R4:S4).
# ---------------------------------------------------------------------
$ws.Range("R4:S4").ClearContents()

# ---------------------------------------------------------------------
# Row 5: drop the leftover Actual/Pass-Fail columns (R5:S5).
# ---------------------------------------------------------------------
$ws.Range("R5:S5").ClearContents()

# ---------------------------------------------------------------------
# Row 6 (new): another Decision_Table / RelatedPartyEnforcedPairs search,
# this time for Param1 "BRO"; Q6 uses the text-numfmt quotePrefix style
# (s=1), matching row 5's Q style.
# ---------------------------------------------------------------------
$ws.Range("A5:E5").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").Value = "Decision_Table"
$ws.Range("B6").Value = "PegaFS-Data-RelCodes"
$ws.Range("C6").Value = "RelatedPartyEnforcedPairs"
$ws.Range("D6").Value = "PegaFS:08-06-01"
$ws.Range("E6").Value = "BRO"

$ws.Range("Q5").Copy()
$ws.Range("Q6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("Q6").Value = "'true"

# ---------------------------------------------------------------------
# Row 7 (new): SLA / xyz search row, same shape as row 4 (quote-prefixed
# "0" fill across Param1-12, D7 left/top-aligned, Q7 general+quotePrefix).
# ---------------------------------------------------------------------
$ws.Range("A4:Q4").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A7").Value = "SLA"
$ws.Range("B7").Value = "PegaFS-Data-RelCodes"
$ws.Range("C7").Value = "xyz"
$ws.Range("D7").Value = "BFSCS:01-01-01"
$ws.Range("E7:P7").Value = "'0"
$ws.Range("Q7").Value = "'true"

# ---------------------------------------------------------------------
# Row 8 (new): SLA / TestSLA4 search row - same shape as row 4/7, except
# Param3/Param4 (G8:H8) are "1" instead of "0".
# ---------------------------------------------------------------------
$ws.Range("A4:Q4").Copy()
$ws.Range("A8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "SLA"
$ws.Range("B8").Value = "PegaFS-Data-RelCodes"
$ws.Range("C8").Value = "TestSLA4"
$ws.Range("D8").Value = "BFSCS:01-01-01"
$ws.Range("E8:P8").Value = "'0"
$ws.Range("G8:H8").Value = "'1"
$ws.Range("Q8").Value = "'true"

# ---------------------------------------------------------------------
# Sheet view: no more frozen/scrolled-to column C, selection lands on A4.
# ---------------------------------------------------------------------
[void]$ws.Range("A1").Select()
[void]$ws.Range("A4").Select()

$excel.CutCopyMode = $false
